# "update ops for external url"
#
# - resources!C2 gains the DEQM MeasureReport profile note (previously empty).
# - ops!B2 / ops!B3 (the submit-data / data-requirements OperationDefinition
#   URLs) are repointed from the old davinci-deqm STU3 URLs to the new
#   external http://hl7.org/fhir/OperationDefinition/Measure-* URLs.
# - Trailing view-state (active sheet / selected cell) is left where the
#   author's edits ended up: meta -> B10, resources -> B9, ops -> B3 (active).

$wb = $excel.ActiveWorkbook

# --- "resources": document that the producer client SHOULD support the ---
# --- DEQM MeasureReport profiles (fill in the previously-blank C2 cell) ---
$wsRes = $wb.Worksheets.Item("resources")
$wsRes.Range("C2").Value = "The DaVinci DEQM Producer Client **SHOULD** be capable of supporting the DEQM MeasureReport Profiles and all the DEQM and QI Core Profiles they reference."

# Setting .Value on a previously-empty cell can swap its style out from
# under it (dropping the quotePrefix/wrap formatting it shared with C3) -
# copy the formatting back from the still-untouched sibling cell C3 so the
# style index is preserved exactly.
$wsRes.Range("C3").Copy()
$wsRes.Range("C2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$wsRes.Rows.Item(2).RowHeight = 45

# --- "ops": point the two operation definitions at the new external URLs ---
$wsOps = $wb.Worksheets.Item("ops")
$wsOps.Range("B2").Value = "http://hl7.org/fhir/OperationDefinition/Measure-submit-data"
$wsOps.Range("B3").Value = "http://hl7.org/fhir/OperationDefinition/Measure-data-requirements"

# --- leave the cursor / active-sheet state where the author left it ---
$wsMeta = $wb.Worksheets.Item("meta")
$wsMeta.Activate() | Out-Null
$wsMeta.Range("B10").Select() | Out-Null

$wsRes.Activate() | Out-Null
$wsRes.Range("B9").Select() | Out-Null

$wsOps.Activate() | Out-Null
$wsOps.Range("B3").Select() | Out-Null
